# Update odds values in Jogos_da_Semana_FlashScore_2024-10-03.xlsx
# per the commit diff (Atualizando o arquivo XLSX).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("L3").Value = 6.5  # was 6
$ws.Range("M3").Value = 1.06  # was 1.07
$ws.Range("N3").Value = 10  # was 8.5
$ws.Range("U3").Value = 2.25  # was 2.2
$ws.Range("V3").Value = 1.57  # was 1.62
$ws.Range("AE3").Value = 23  # was 21
$ws.Range("AG3").Value = 11  # was 12
$ws.Range("AZ3").Value = 151  # was 126

# Row 4
$ws.Range("G4").Value = 1.62  # was 1.83
$ws.Range("H4").Value = 3.3  # was 3
$ws.Range("I4").Value = 7  # was 5.25
$ws.Range("J4").Value = 2.38  # was 2.63
$ws.Range("K4").Value = 1.95  # was 1.91
$ws.Range("L4").Value = 7  # was 6
$ws.Range("Q4").Value = 2.7  # was 2.88
$ws.Range("R4").Value = 1.44  # was 1.4
$ws.Range("U4").Value = 2.63  # was 2.5
$ws.Range("V4").Value = 1.44  # was 1.5
$ws.Range("W4").Value = 4.5  # was 4.75
$ws.Range("X4").Value = 6  # was 7
$ws.Range("Z4").Value = 11  # was 15
$ws.Range("AA4").Value = 19  # was 21
$ws.Range("AD4").Value = 7  # was 6.5
$ws.Range("AE4").Value = 29  # was 23
$ws.Range("AF4").Value = 126  # was 101
$ws.Range("AG4").Value = 11  # was 9.5
$ws.Range("AH4").Value = 29  # was 23
$ws.Range("AI4").Value = 23  # was 19
$ws.Range("AJ4").Value = 81  # was 51
$ws.Range("AK4").Value = 67  # was 51
$ws.Range("AL4").Value = 81  # was 67
$ws.Range("AN4").Value = 3.25  # was 3.6
$ws.Range("AO4").Value = 9  # was 11
$ws.Range("AQ4").Value = 34  # was 41
$ws.Range("AR4").Value = 67  # was 81
$ws.Range("AU4").Value = 11  # was 10
$ws.Range("AW4").Value = 7.5  # was 6.5
$ws.Range("AX4").Value = 41  # was 34
$ws.Range("AY4").Value = 51  # was 41
$ws.Range("AZ4").Value = 201  # was 126
$ws.Range("BA4").Value = 301  # was 201

# Row 5
$ws.Range("G5").Value = 1.67  # was 1.57
$ws.Range("H5").Value = 3.7  # was 3.8
$ws.Range("I5").Value = 4.2  # was 4.75
$ws.Range("J5").Value = 2.3  # was 2.2
$ws.Range("L5").Value = 4.75  # was 5.5
$ws.Range("Q5").Value = 1.95  # was 1.93
$ws.Range("R5").Value = 1.85  # was 1.88
$ws.Range("U5").Value = 1.91  # was 2
$ws.Range("V5").Value = 1.8  # was 1.73
$ws.Range("W5").Value = 7  # was 6.5
$ws.Range("X5").Value = 8  # was 7.5
$ws.Range("Z5").Value = 13  # was 12
$ws.Range("AD5").Value = 7.5  # was 8
$ws.Range("AE5").Value = 17  # was 19
$ws.Range("AF5").Value = 51  # was 67
$ws.Range("AH5").Value = 23  # was 26
$ws.Range("AI5").Value = 15  # was 17
$ws.Range("AM5").Value = 900  # was 1000
$ws.Range("AN5").Value = 3.75  # was 3.6
$ws.Range("AO5").Value = 9  # was 8.5
$ws.Range("AQ5").Value = 29  # was 26
$ws.Range("AU5").Value = 8.5  # was 9
$ws.Range("AW5").Value = 6.5  # was 7
$ws.Range("AX5").Value = 26  # was 29
$ws.Range("AZ5").Value = 81  # was 101
$ws.Range("BA5").Value = 101  # was 126

# Row 6
$ws.Range("G6").Value = 4.5  # was 5
$ws.Range("H6").Value = 3.8  # was 4
$ws.Range("I6").Value = 1.62  # was 1.53
$ws.Range("J6").Value = 4.33  # was 4.75
$ws.Range("L6").Value = 2.1  # was 2
$ws.Range("M6").Value = 1.03  # was 1.02
$ws.Range("N6").Value = 10.5  # was 11
$ws.Range("U6").Value = 1.53  # was 1.57
$ws.Range("V6").Value = 2.38  # was 2.25
$ws.Range("W6").Value = 19  # was 21
$ws.Range("X6").Value = 29  # was 34
$ws.Range("Y6").Value = 15  # was 17
$ws.Range("AC6").Value = 17  # was 19
$ws.Range("AF6").Value = 34  # was 41
$ws.Range("AG6").Value = 11  # was 10
$ws.Range("AH6").Value = 10  # was 9.5
$ws.Range("AJ6").Value = 13  # was 12
$ws.Range("AM6").Value = 101  # was 126
$ws.Range("AN6").Value = 6.5  # was 7
$ws.Range("AO6").Value = 21  # was 23
$ws.Range("AP6").Value = 23  # was 26
$ws.Range("AR6").Value = 67  # was 81
$ws.Range("AX6").Value = 8  # was 7.5
$ws.Range("AZ6").Value = 23  # was 21
$ws.Range("BA6").Value = 41  # was 34

# Row 7
$ws.Range("M7").Value = 1.02  # was 1.04
$ws.Range("N7").Value = 11  # was 9
$ws.Range("O7").Value = 1.25  # was 1.22
$ws.Range("P7").Value = 3.75  # was 4
$ws.Range("Q7").Value = 1.8  # was 1.75
$ws.Range("R7").Value = 2  # was 2.05
$ws.Range("S7").Value = 1.36  # was 1.33
$ws.Range("T7").Value = 3  # was 3.25
$ws.Range("U7").Value = 1.67  # was 1.62
$ws.Range("V7").Value = 2.1  # was 2.2
$ws.Range("W7").Value = 9  # was 9.5
$ws.Range("AB7").Value = 26  # was 23
$ws.Range("AC7").Value = 11  # was 12
$ws.Range("AG7").Value = 11  # was 12
$ws.Range("AT7").Value = 3  # was 3.25
$ws.Range("AX7").Value = 17  # was 15

# Row 8
$ws.Range("M8").Value = 1.03  # was 1.04
$ws.Range("N8").Value = 15  # was 13
$ws.Range("Q8").Value = 1.67  # was 1.7
$ws.Range("R8").Value = 2.15  # was 2.1
$ws.Range("W8").Value = 12  # was 11
$ws.Range("AC8").Value = 15  # was 13
$ws.Range("AM8").Value = 126  # was 151
$ws.Range("AQ8").Value = 41  # was 51
$ws.Range("BC8").Value = 401  # was 451
